# Update data: 11 July 2020
# Append the newest observation (2020-06-01, Excel serial 43983) to both the
# "Canada" sheet (national figures) and the "Province" sheet (per-province
# breakdown), matching the monthly append pattern already used in the file.

$wb = $excel.ActiveWorkbook

$dateFmt = "d-mmm-yy"
$newDate = 43983

# ---------------------------------------------------------------------
# Sheet "Canada": one new row (row 7) with the national figure.
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A7").Value = $newDate
$wsCanada.Range("A7").NumberFormat = $dateFmt

$wsCanada.Range("B7").Value = "Canada"
$wsCanada.Range("B7").NumberFormat = $dateFmt

$wsCanada.Range("C7").Value = 118.4
$wsCanada.Range("D7").Value = 1122.9000000000001

[void]$wsCanada.Range("C8").Select()

# ---------------------------------------------------------------------
# Sheet "Province": ten new rows (52-61), one per province, in the same
# order used for every earlier month in the sheet.
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 52; Name = "Newfoundland & Labrador"; C = 18.6;  D = 33.9 },
    @{ Row = 53; Name = "Prince Edward Island";     C = 69.2;  D = 7.8 },
    @{ Row = 54; Name = "Nova Scotia";              C = 91.1;  D = 33.799999999999997 },
    @{ Row = 55; Name = "New Brunswick";            C = 25.6;  D = 30.8 },
    @{ Row = 56; Name = "Quebec";                   C = 116.4; D = 224.3 },
    @{ Row = 57; Name = "Ontario";                  C = 119;   D = 430.8 },
    @{ Row = 58; Name = "Manitoba";                 C = 77.2;  D = 39 },
    @{ Row = 59; Name = "Saskatchewan";              C = 121.8; D = 31.7 },
    @{ Row = 60; Name = "Alberta";                  C = 127.9; D = 167.8 },
    @{ Row = 61; Name = "British Columbia";         C = 180.9; D = 123 }
)

foreach ($entry in $provinceRows) {
    $r = $entry.Row

    $wsProvince.Cells.Item($r, 1).Value = $newDate
    $wsProvince.Cells.Item($r, 1).NumberFormat = $dateFmt

    $wsProvince.Cells.Item($r, 2).Value = $entry.Name
    if ($r -eq 52) {
        # Only the first row of each monthly block carries the date style
        # on column B too (matches the pattern of every earlier block).
        $wsProvince.Cells.Item($r, 2).NumberFormat = $dateFmt
    }

    $wsProvince.Cells.Item($r, 3).Value = $entry.C
    $wsProvince.Cells.Item($r, 4).Value = $entry.D
}

[void]$wsProvince.Range("A61").Select()
